# RIR.docx: add the missing "Condition" section (a "Condition:" label
# paragraph followed by a 2-column "Condtion Value" / "$Condition_Value"
# table) right after the existing "Quality of Service" table, mirroring
# that table's own layout/formatting. The trailing bookmark paragraph
# (_GoBack) that closes the document body is left untouched.

$d = $word.ActiveDocument

# Guard against re-running this edit on a document that already has it.
if ($d.Content.Text -notmatch [regex]::Escape("Condition:")) {

    # The existing "Quality of Service" table is the last (only) table
    # currently in the body; insert immediately after it.
    $qosTable = $d.Tables.Item($d.Tables.Count)

    $insertPoint = $d.Range($qosTable.Range.End, $qosTable.Range.End)

    $wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

    $fragment = (
        '<w:p ' + $wNs + '/>' +
        '<w:p ' + $wNs + '><w:r><w:t>Condition:</w:t></w:r></w:p>' +
        '<w:tbl ' + $wNs + '>' +
            '<w:tblPr>' +
                '<w:tblW w:w="5000" w:type="pct"/>' +
                '<w:tblLayout w:type="fixed"/>' +
                '<w:tblLook w:val="00A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/>' +
            '</w:tblPr>' +
            '<w:tblGrid>' +
                '<w:gridCol w:w="3231"/>' +
                '<w:gridCol w:w="6113"/>' +
            '</w:tblGrid>' +
            '<w:tr>' +
                '<w:tc>' +
                    '<w:tcPr>' +
                        '<w:tcW w:w="1729" w:type="pct"/>' +
                        '<w:tcBorders>' +
                            '<w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                            '<w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                            '<w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                            '<w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                        '</w:tcBorders>' +
                        '<w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>' +
                    '</w:tcPr>' +
                    '<w:p>' +
                        '<w:pPr>' +
                            '<w:spacing w:line="276" w:lineRule="auto"/>' +
                            '<w:rPr>' +
                                '<w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/>' +
                                '<w:b/>' +
                                '<w:sz w:val="22"/>' +
                                '<w:szCs w:val="22"/>' +
                            '</w:rPr>' +
                        '</w:pPr>' +
                        '<w:r>' +
                            '<w:rPr>' +
                                '<w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/>' +
                                '<w:b/>' +
                                '<w:sz w:val="22"/>' +
                                '<w:szCs w:val="22"/>' +
                            '</w:rPr>' +
                            '<w:t>Condtion Value</w:t>' +
                        '</w:r>' +
                    '</w:p>' +
                '</w:tc>' +
                '<w:tc>' +
                    '<w:tcPr>' +
                        '<w:tcW w:w="3271" w:type="pct"/>' +
                        '<w:tcBorders>' +
                            '<w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                            '<w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                            '<w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                            '<w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>' +
                        '</w:tcBorders>' +
                        '<w:shd w:val="clear" w:color="auto" w:fill="auto"/>' +
                    '</w:tcPr>' +
                    '<w:p>' +
                        '<w:pPr>' +
                            '<w:rPr>' +
                                '<w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/>' +
                                '<w:sz w:val="22"/>' +
                                '<w:szCs w:val="22"/>' +
                                '<w:lang w:eastAsia="zh-CN"/>' +
                            '</w:rPr>' +
                        '</w:pPr>' +
                        '<w:r><w:t>$Condition_Value</w:t></w:r>' +
                    '</w:p>' +
                '</w:tc>' +
            '</w:tr>' +
        '</w:tbl>'
    )

    $null = $insertPoint.InsertXML($fragment)

    Write-Output "Inserted Condition label + table after the Quality of Service table."
} else {
    Write-Output "Condition section already present; no changes made."
}
